$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.012.03'
$ws.Range('E2').Value = '  +16.55%  '
$ws.Range('D3').Value = '1.664.18'
$ws.Range('E3').Value = '  +12.50%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.71%  '
$ws.Range('D5').Value = '''307.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.76%  '
$ws.Range('D6').Value = '''0.9963'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.68%  '
$ws.Range('D7').Value = '''0.3718'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.90%  '
$ws.Range('D8').Value = '''0.3446'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.94%  '
$ws.Range('D9').Value = '''47.73'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +21.08%  '
$ws.Range('D10').Value = '''1.168'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.87%  '
$ws.Range('D11').Value = '''0.07243'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +8.77%  '
$ws.Range('D12').Value = '''0.9978'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = '''20.48'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +12.64%  '
$ws.Range('D14').Value = '''6.022'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +10.12%  '
$ws.Range('E15').Value = '  +8.81%  '
$ws.Range('D16').Value = '1.659.19'
$ws.Range('E16').Value = '  +12.41%  '
$ws.Range('D17').Value = '''0.00001095'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +7.86%  '
$ws.Range('D18').Value = '''0.9961'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('D19').Value = '''0.06714'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.97%  '
$ws.Range('D20').Value = '''81.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +18.33%  '
$ws.Range('D21').Value = '''16.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +13.06%  '
$ws.Range('D22').Value = '''6.114'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +11.48%  '
$ws.Range('D23').Value = '''12.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.05%  '
$ws.Range('D24').Value = '23.966.59'
$ws.Range('E24').Value = '  +16.34%  '
$ws.Range('D25').Value = '''2.385'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.98%  '
$ws.Range('D26').Value = '''3.385'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -8.01%  '
$ws.Range('D27').Value = '''2.667'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +28.35%  '
$ws.Range('D28').Value = '''152.11'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.48%  '
$ws.Range('D29').Value = '''19.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +13.54%  '
$ws.Range('D30').Value = '1.841.07'
$ws.Range('E30').Value = '  +12.50%  '
$ws.Range('D31').Value = '''126.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.45%  '
$ws.Range('D32').Value = '''6.357'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +28.55%  '
$ws.Range('D33').Value = '''4.094'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.47%  '
$ws.Range('D34').Value = '''0.9763'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +22.55%  '
$ws.Range('E35').Value = '  +21.52%  '
$ws.Range('D36').Value = '''0.08388'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.81%  '
$ws.Range('D37').Value = '''12.32'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +20.09%  '
$ws.Range('D38').Value = '''8.958'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +21.51%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '''5.320'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.69%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '''0.06352'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.78%  '
$ws.Range('D41').Value = '''1.291'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.70%  '
$ws.Range('D42').Value = '''0.02316'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +14.81%  '
$ws.Range('D43').Value = '''0.2077'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +12.63%  '
$ws.Range('D44').Value = '''0.6098'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +16.67%  '
$ws.Range('D45').Value = '''0.9960'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.73%  '
$ws.Range('D46').Value = '''3.809'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.37%  '
$ws.Range('D47').Value = '''13.20'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.64%  '
$ws.Range('D48').Value = '''0.5941'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +14.87%  '
$ws.Range('D49').Value = '''126.89'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.73%  '
$ws.Range('D50').Value = '''2.001'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.35%  '
$ws.Range('D51').Value = '''0.07083'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.51%  '
